$d = $word.ActiveDocument
$t = $d.Tables(1)

# wdAlignParagraphCenter
$wdCenter = 1

function Set-CellCentered($cell, $text) {
    $cell.Range.Text = $text
    $cell.Range.ParagraphFormat.Alignment = $wdCenter
}

# ---- Row "2" executors (table row 3) : dataset 20M ----
$row = $t.Rows(3)
$row.Cells(1).Range.Paragraphs(1).Range.Text = "2 (20M)"
Set-CellCentered $row.Cells(2) "304956"
Set-CellCentered $row.Cells(3) "603"
Set-CellCentered $row.Cells(4) "171505"
Set-CellCentered $row.Cells(5) "12271"

# ---- Row "4" executors (table row 4) : dataset 20M ----
$row = $t.Rows(4)
$row.Cells(1).Range.Paragraphs(1).Range.Text = "4 (20M)"
Set-CellCentered $row.Cells(2) "157964"
Set-CellCentered $row.Cells(3) "390"
Set-CellCentered $row.Cells(4) "84303"
Set-CellCentered $row.Cells(5) "6501"

# ---- Row "8" executors (table row 5) : dataset 50M ----
$row = $t.Rows(5)
$row.Cells(1).Range.Paragraphs(1).Range.Text = "8 (50M)"
Set-CellCentered $row.Cells(2) "190337"
Set-CellCentered $row.Cells(3) "231"
Set-CellCentered $row.Cells(4) "112348"
Set-CellCentered $row.Cells(5) "7210"

# ---- Row "16" executors (table row 6) : dataset 100M ----
$row = $t.Rows(6)
$row.Cells(1).Range.Paragraphs(1).Range.Text = "16 (100M)"
Set-CellCentered $row.Cells(2) "200158"
Set-CellCentered $row.Cells(3) "191"
Set-CellCentered $row.Cells(4) "113742"
Set-CellCentered $row.Cells(5) "7151"

# ---- Add trailing summary paragraph after the last (empty) list paragraph ----
$newPara = $d.Paragraphs.Add()
$newPara.Range.Text = "2 exec 20 Milioni: R1 298285 R2 602 R3 170077 MR 11951"
